$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14..90 down to 15..91
$ws.Rows.Item(14).Insert(-4121)

# Populate the new row 14 with a new weekly data record.
# Columns that stay identical to the (old) row 14 are simply re-written with
# the same values; columns that actually change carry the new values from
# the diff (Fecha, Precio minimo/maximo/promedio, Origen, Precio $/Kg).
$ws.Cells.Item(14, 1).Value = 7
$ws.Cells.Item(14, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value = "Ñuble"
$ws.Cells.Item(14, 4).Value = 44901
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = 100112022
$ws.Cells.Item(14, 7).Value = "Arveja Verde"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 80
$ws.Cells.Item(14, 11).Value = 20000
$ws.Cells.Item(14, 12).Value = 22000
$ws.Cells.Item(14, 13).Value = 21000
$ws.Cells.Item(14, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 840
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
